# Add a new "jammer" sheet (data/init_location.xlsx upload) with initial
# jammer location data, mirroring the layout used by the other location
# sheets (index, x, y, z).

$wb = $excel.ActiveWorkbook

# --- Update selections on a couple of the existing sheets -----------------
$wsAttacker = $wb.Worksheets.Item("attacker")
$wsAttacker.Range("A2:D3").Select()

$wsRIS = $wb.Worksheets.Item("RIS")
$wsRIS.Range("A1:D1").Select()

# --- Add the new "jammer" sheet, placed after the last existing sheet -----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsJammer = $wb.Worksheets.Add($null, $lastSheet)
$wsJammer.Name = "jammer"

# Header row
$wsJammer.Range("A1").Value = "index"
$wsJammer.Range("B1").Value = "x"
$wsJammer.Range("C1").Value = "y"
$wsJammer.Range("D1").Value = "z"

# Data rows
$data = @(
    @(0, -4, 47, 0.0001),
    @(1, 20, 25, 0.0001),
    @(2, 20, 47, 0.0001),
    @(3, 20, 25, 0.0001),
    @(4, 4, 47, 0.0001)
)

$r = 2
foreach ($row in $data) {
    $wsJammer.Cells.Item($r, 1).Value = $row[0]
    $wsJammer.Cells.Item($r, 2).Value = $row[1]
    $wsJammer.Cells.Item($r, 3).Value = $row[2]
    $wsJammer.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsJammer.Range("F17").Select()

$wsJammer.Activate()
